$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "59.964.42"
$c.ClearFormats()

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +3.73%  "
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.419.44"
$c.ClearFormats()

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +3.28%  "
$c.ClearFormats()

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  -0.01%  "
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "551.84"
$c.ClearFormats()

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +2.32%  "
$c.ClearFormats()

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "137.87"
$c.ClearFormats()

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +2.77%  "
$c.ClearFormats()

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -0.06%  "
$c.ClearFormats()

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.576"
$c.ClearFormats()

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +1.18%  "
$c.ClearFormats()

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +3.99%  "
$c.ClearFormats()

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "5.76"
$c.ClearFormats()

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +4.52%  "
$c.ClearFormats()

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.360"
$c.ClearFormats()

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +1.11%  "
$c.ClearFormats()

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -1.92%  "
$c.ClearFormats()

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "24.90"
$c.ClearFormats()

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +5.00%  "
$c.ClearFormats()

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.848.52"
$c.ClearFormats()

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  +3.15%  "
$c.ClearFormats()

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "59.907.47"
$c.ClearFormats()

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +3.73%  "
$c.ClearFormats()

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000138"
$c.ClearFormats()

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +3.02%  "
$c.ClearFormats()

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.429.58"
$c.ClearFormats()

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +3.58%  "
$c.ClearFormats()

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "11.34"
$c.ClearFormats()

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +6.30%  "
$c.ClearFormats()

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +2.43%  "
$c.ClearFormats()

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +1.18%  "
$c.ClearFormats()

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +1.16%  "
$c.ClearFormats()

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +0.19%  "
$c.ClearFormats()

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "65.04"
$c.ClearFormats()

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +3.57%  "
$c.ClearFormats()

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +4.07%  "
$c.ClearFormats()

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "8.58"
$c.ClearFormats()

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +3.19%  "
$c.ClearFormats()

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +0.41%  "
$c.ClearFormats()

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  +0.22%  "
$c.ClearFormats()

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.0₃0782"
$c.ClearFormats()

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +6.48%  "
$c.ClearFormats()

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +1.47%  "
$c.ClearFormats()

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.26"
$c.ClearFormats()

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  +2.57%  "
$c.ClearFormats()

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "169.15"
$c.ClearFormats()

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -0.49%  "
$c.ClearFormats()

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +2.83%  "
$c.ClearFormats()

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "18.73"
$c.ClearFormats()

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +2.25%  "
$c.ClearFormats()

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +0.02%  "
$c.ClearFormats()

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +5.47%  "
$c.ClearFormats()

$c = $ws.Range("B36")
$c.NumberFormat = "@"
$c.Value = "FirstDigitalUSD"
$c.ClearFormats()

$c = $ws.Range("C36")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +0.09%  "
$c.ClearFormats()

$c = $ws.Range("B37")
$c.NumberFormat = "@"
$c.Value = "NEARProtocol"
$c.ClearFormats()

$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c.ClearFormats()

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.21"
$c.ClearFormats()

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +0.98%  "
$c.ClearFormats()

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  +0.40%  "
$c.ClearFormats()

$c = $ws.Range("B39")
$c.NumberFormat = "@"
$c.Value = "OKB"
$c.ClearFormats()

$c = $ws.Range("C39")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c.ClearFormats()

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "39.55"
$c.ClearFormats()

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  +1.24%  "
$c.ClearFormats()

$c = $ws.Range("B40")
$c.NumberFormat = "@"
$c.Value = "Bittensor"
$c.ClearFormats()

$c = $ws.Range("C40")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c.ClearFormats()

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "319.82"
$c.ClearFormats()

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +11.34%  "
$c.ClearFormats()

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  +10.49%  "
$c.ClearFormats()

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.68"
$c.ClearFormats()

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +1.28%  "
$c.ClearFormats()

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "139.76"
$c.ClearFormats()

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -1.15%  "
$c.ClearFormats()

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +1.42%  "
$c.ClearFormats()

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0520"
$c.ClearFormats()

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +1.84%  "
$c.ClearFormats()

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "19.56"
$c.ClearFormats()

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +2.50%  "
$c.ClearFormats()

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.411"
$c.ClearFormats()

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  +7.94%  "
$c.ClearFormats()

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +1.51%  "
$c.ClearFormats()

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0227"
$c.ClearFormats()

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +2.19%  "
$c.ClearFormats()

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +1.57%  "
$c.ClearFormats()

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -0.13%  "
$c.ClearFormats()

